$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.579.02"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "2.666.21"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.18"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.86"
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("E8").Value = "  +4.39%  "
$ws.Range("E9").Value = "  +2.20%  "
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("E11").Value = "  -2.94%  "
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.99"
$ws.Range("E13").Value = "  -3.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000197"
$ws.Range("E14").Value = "  -4.29%  "
$ws.Range("D15").Value = "3.146.54"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").Value = "65.488.23"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").Value = "2.661.89"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.64"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("E19").Value = "  -1.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.48"
$ws.Range("E20").Value = "  -2.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350.72"
$ws.Range("E21").Value = "  -1.65%  "
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.11"
$ws.Range("E23").Value = "  -2.45%  "
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("E25").Value = "  -2.37%  "
$ws.Range("E26").Value = "  +2.95%  "
$ws.Range("E27").Value = "  -3.38%  "
$ws.Range("E28").Value = "  -3.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.97"
$ws.Range("E29").Value = "  -3.22%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.12"
$ws.Range("E31").Value = "  -3.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "530.30"
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("E33").Value = "  +0.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.43"
$ws.Range("E34").Value = "  -2.53%  "
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.421"
$ws.Range("E36").Value = "  -2.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.56"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "156.85"
$ws.Range("E39").Value = "  -2.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.92"
$ws.Range("E40").Value = "  -3.42%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "162.36"
$ws.Range("E42").Value = "  -2.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.08"
$ws.Range("E43").Value = "  -1.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.32"
$ws.Range("E44").Value = "  +2.01%  "
$ws.Range("E45").Value = "  -3.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.54"
$ws.Range("E46").Value = "  -4.28%  "
$ws.Range("E47").Value = "  -2.41%  "
$ws.Range("E48").Value = "  -3.05%  "
$ws.Range("D49").Value = "0.0₆0252"
$ws.Range("E49").Value = "  +6.23%  "
$ws.Range("E50").Value = "  -0.91%  "
$ws.Range("E51").Value = "  -4.03%  "
